$d = $word.ActiveDocument

$replacements = @(
    @("81÷3=", "38÷6="),
    @("19÷5=", "84÷4="),
    @("65÷4=", "23÷8="),
    @("93÷8=", "99÷2="),
    @("79÷8=", "65÷7="),
    @("76÷7=", "52÷7="),
    @("98÷4=", "39÷9="),
    @("12÷6=", "49÷6="),
    @("93÷3=", "67÷9="),
    @("47÷4=", "73÷8="),
    @("49÷7=", "73÷4="),
    @("29÷2=", "94÷5="),
    @("75÷8=", "38÷6="),
    @("39÷8=", "32÷5="),
    @("28÷7=", "36÷5="),
    @("41÷7=", "13÷7="),
    @("29÷9=", "35÷5="),
    @("88÷2=", "99÷6="),
    @("67÷4=", "65÷6="),
    @("44÷3=", "39÷7="),
    @("69÷2=", "21÷2="),
    @("55÷5=", "72÷5="),
    @("16÷4=", "16÷5="),
    @("87÷9=", "62÷7="),
    @("86÷6=", "94÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
